$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text cells keep their original text type
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '69.896.82'
$ws.Range("E2").Value = '  +5.31%  '
$ws.Range("D3").Value = '3.404.15'
$ws.Range("E3").Value = '  +10.55%  '
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value = '585.97'
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").Value = '182.58'
$ws.Range("E6").Value = '  +7.50%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '3.415.19'
$ws.Range("E8").Value = '  +11.06%  '
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  +4.45%  '
$ws.Range("D10").Value = '6.62'
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("D11").Value = '0.157'
$ws.Range("E11").Value = '  +4.76%  '
$ws.Range("D12").Value = '0.488'
$ws.Range("E12").Value = '  +3.88%  '
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  +4.89%  '
$ws.Range("D14").Value = '38.16'
$ws.Range("E14").Value = '  +6.64%  '
$ws.Range("D15").Value = '3.945.70'
$ws.Range("E15").Value = '  +9.80%  '
$ws.Range("D16").Value = '69.592.40'
$ws.Range("E16").Value = '  +4.96%  '
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '7.40'
$ws.Range("E18").Value = '  +6.26%  '
$ws.Range("D19").Value = '3.339.32'
$ws.Range("E19").Value = '  +8.38%  '
$ws.Range("D20").Value = '17.10'
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("D21").Value = '503.46'
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("D22").Value = '8.21'
$ws.Range("E22").Value = '  +6.18%  '
$ws.Range("D23").Value = '0.727'
$ws.Range("E23").Value = '  +5.81%  '
$ws.Range("D24").Value = '86.16'
$ws.Range("E24").Value = '  +4.20%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '13.39'
$ws.Range("E25").Value = '  +5.71%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  +8.70%  '
$ws.Range("D27").Value = '10.77'
$ws.Range("E27").Value = '  +5.85%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").Value = '8.25'
$ws.Range("E29").Value = '  +5.26%  '
$ws.Range("D30").Value = '2.48'
$ws.Range("E30").Value = '  +9.46%  '
$ws.Range("D31").Value = '2.69'
$ws.Range("E31").Value = '  +3.17%  '
$ws.Range("D32").Value = '29.84'
$ws.Range("E32").Value = '  +7.94%  '
$ws.Range("D33").Value = '0.0000101'
$ws.Range("E33").Value = '  +12.02%  '
$ws.Range("D34").Value = '0.116'
$ws.Range("E34").Value = '  +4.66%  '
$ws.Range("D35").Value = '0.992'
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("D36").Value = '6.01'
$ws.Range("E36").Value = '  +7.41%  '
$ws.Range("D37").Value = '1.01'
$ws.Range("E37").Value = '  +6.19%  '
$ws.Range("D38").Value = '48.77'
$ws.Range("E38").Value = '  +2.64%  '
$ws.Range("D39").Value = '0.330'
$ws.Range("E39").Value = '  +9.79%  '
$ws.Range("D40").Value = '2.11'
$ws.Range("E40").Value = '  +7.44%  '
$ws.Range("E41").Value = '  +4.79%  '
$ws.Range("D42").Value = '50.16'
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("D43").Value = '8.65'
$ws.Range("E43").Value = '  +4.48%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '413.16'
$ws.Range("E44").Value = '  +13.06%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  +11.57%  '
$ws.Range("D46").Value = '2.925.50'
$ws.Range("E46").Value = '  +4.02%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0359'
$ws.Range("E47").Value = '  +3.93%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '27.50'
$ws.Range("E48").Value = '  +12.75%  '
$ws.Range("D49").Value = '135.83'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D51").Value = '2.43'
$ws.Range("E51").Value = '  +12.89%  '

# Restore default (unstyled) cell style, same as original cells
$ws.Range("D2:E51").Style = "Normal"
